# Clean up vaccine category labels across all worksheets:
#  - remove footnote markers like " [1]", "[2]", "[5]", etc.
#  - collapse embedded line breaks into a single space
# This also naturally fixes the mis-typed "Afluria\nQuadrivalent" entry so
# that it reads the same as the other "Afluria Quadrivalent" cells, letting
# Excel drop the now-unused shared string automatically on save.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value()

            if ($val -ne $null -and $val.GetType().Name -eq "String") {
                $newVal = $val -replace '\[\d+\]', ''
                $newVal = $newVal -replace "`n", ' '

                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
